$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.98
$ws.Range("G2").Value = 2.52
$ws.Range("H2").Value = 3.35
$ws.Range("P2").Value = 1.79

# Row 3
$ws.Range("G3").Value = 2.68
$ws.Range("H3").Value = 3.15
$ws.Range("I3").Value = 3.95
$ws.Range("J3").Value = 3.25
$ws.Range("K3").Value = 4.1
$ws.Range("N3").Value = 3.5
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.84
$ws.Range("Q3").Value = 1.83

# Row 4
$ws.Range("P4").Value = 2.4
$ws.Range("Q4").Value = 1.41

# Row 6
$ws.Range("T6").Value = 1.85
$ws.Range("V6").Value = 1.37

# Row 7
$ws.Range("G7").Value = 1.56
$ws.Range("P7").Value = 2.38
$ws.Range("Q7").Value = 1.49

# Row 8
$ws.Range("F8").Value = 1.38
$ws.Range("G8").Value = 1.53
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 12
$ws.Range("J8").Value = 4.4
$ws.Range("K8").Value = 6.6
$ws.Range("N8").Value = 4.7
$ws.Range("P8").Value = 2.28
$ws.Range("Q8").Value = 1.54

# Row 9
$ws.Range("G9").Value = 6.4
$ws.Range("H9").Value = 1.68
$ws.Range("J9").Value = 3.65

# Row 10
$ws.Range("Q10").Value = 1.42

# Row 11
$ws.Range("F11").Value = 1.49
$ws.Range("G11").Value = 1.55
$ws.Range("H11").Value = 8.4
$ws.Range("P11").Value = 1.68
$ws.Range("Q11").Value = 2.22

# Row 15
$ws.Range("F15").Value = 1.54
$ws.Range("G15").Value = 1.79
$ws.Range("H15").Value = 6.4
$ws.Range("I15").Value = 8.4
$ws.Range("J15").Value = 3.4
$ws.Range("K15").Value = 5
$ws.Range("N15").Value = 3.2
$ws.Range("O15").Value = 1.34
$ws.Range("P15").Value = 1.87
$ws.Range("Q15").Value = 1.99

# Row 16
$ws.Range("N16").Value = 3.1
$ws.Range("O16").Value = 1.38

# Row 17
$ws.Range("J17").Value = 2.98

# Row 18
$ws.Range("N18").Value = 3.15

# Row 19
$ws.Range("F19").Value = 1.93
$ws.Range("G19").Value = 2.02
$ws.Range("H19").Value = 4.1
$ws.Range("T19").Value = 1.71
$ws.Range("X19").Value = 18
$ws.Range("Z19").Value = 36
$ws.Range("AB19").Value = 10.5
$ws.Range("AG19").Value = 11
$ws.Range("AH19").Value = 19
$ws.Range("AK19").Value = 21

# Row 20
$ws.Range("N20").Value = 3.7
$ws.Range("U20").Value = 2.16
$ws.Range("AF20").Value = 17.5
$ws.Range("AL20").Value = 46

# Row 21
$ws.Range("R21").Value = 1.63
$ws.Range("S21").Value = 2.5
$ws.Range("U21").Value = 2.54
$ws.Range("Z21").Value = 36
$ws.Range("AJ21").Value = 22

# Row 22
$ws.Range("I22").Value = 24
$ws.Range("K22").Value = 8.6
$ws.Range("P22").Value = 2.36
